# Updated symbol list on Sun Feb  5 05:43:18 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) text values for
# the crypto rows that changed. Values are stored as plain text in the
# workbook (not numbers/percentages), so each new value is written with a
# leading apostrophe to force Excel to keep it as text and avoid any
# numeric/percentage auto-conversion.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'332.22"
$ws.Cells.Item(2, 5).Value = "'1.00%"
$ws.Cells.Item(3, 4).Value = "'44.12"
$ws.Cells.Item(3, 5).Value = "'7.49%"
$ws.Cells.Item(4, 4).Value = "'5.768"
$ws.Cells.Item(4, 5).Value = "'2.98%"
$ws.Cells.Item(5, 4).Value = "'0.08339"
$ws.Cells.Item(5, 5).Value = "'2.13%"
$ws.Cells.Item(6, 4).Value = "'8.800"
$ws.Cells.Item(7, 4).Value = "'4.504"
$ws.Cells.Item(7, 5).Value = "'-0.81%"
$ws.Cells.Item(8, 5).Value = "'-3.81%"
$ws.Cells.Item(10, 4).Value = "'0.9321"
$ws.Cells.Item(10, 5).Value = "'1.42%"
$ws.Cells.Item(11, 4).Value = "'0.1244"
$ws.Cells.Item(11, 5).Value = "'-1.23%"
$ws.Cells.Item(12, 5).Value = "'-0.18%"
$ws.Cells.Item(13, 4).Value = "'0.09494"
$ws.Cells.Item(13, 5).Value = "'1.89%"
$ws.Cells.Item(14, 4).Value = "'0.03958"
$ws.Cells.Item(14, 5).Value = "'5.76%"
$ws.Cells.Item(15, 4).Value = "'0.1065"
$ws.Cells.Item(15, 5).Value = "'0.83%"
$ws.Cells.Item(16, 4).Value = "'0.001302"
$ws.Cells.Item(16, 5).Value = "'0.76%"
$ws.Cells.Item(17, 4).Value = "'0.005927"
$ws.Cells.Item(17, 5).Value = "'-5.36%"
$ws.Cells.Item(18, 4).Value = "'3.504"
$ws.Cells.Item(18, 5).Value = "'1.88%"
$ws.Cells.Item(20, 4).Value = "'9.074"
$ws.Cells.Item(20, 5).Value = "'8.77%"
$ws.Cells.Item(21, 4).Value = "'0.1362"
$ws.Cells.Item(21, 5).Value = "'-2.29%"
$ws.Cells.Item(23, 4).Value = "'0.04416"
$ws.Cells.Item(23, 5).Value = "'-0.05%"
$ws.Cells.Item(24, 4).Value = "'0.001257"
$ws.Cells.Item(24, 5).Value = "'-0.16%"
$ws.Cells.Item(25, 4).Value = "'0.004386"
$ws.Cells.Item(25, 5).Value = "'1.80%"
$ws.Cells.Item(26, 5).Value = "'0.76%"
$ws.Cells.Item(27, 4).Value = "'0.0003993"
$ws.Cells.Item(27, 5).Value = "'-0.01%"
$ws.Cells.Item(39, 4).Value = "'0.02833"
$ws.Cells.Item(39, 5).Value = "'3.23%"
$ws.Cells.Item(40, 4).Value = "'0.05778"
$ws.Cells.Item(40, 5).Value = "'6.94%"
$ws.Cells.Item(41, 4).Value = "'0.007924"
$ws.Cells.Item(41, 5).Value = "'3.34%"
$ws.Cells.Item(42, 4).Value = "'0.1426"
$ws.Cells.Item(42, 5).Value = "'0.80%"
$ws.Cells.Item(43, 4).Value = "'0.009086"
$ws.Cells.Item(43, 5).Value = "'-0.77%"
$ws.Cells.Item(44, 5).Value = "'-1.49%"
$ws.Cells.Item(45, 4).Value = "'0.01017"
$ws.Cells.Item(45, 5).Value = "'-10.32%"
$ws.Cells.Item(46, 4).Value = "'0.00007274"
$ws.Cells.Item(46, 5).Value = "'5.58%"
$ws.Cells.Item(47, 5).Value = "'-0.11%"
$ws.Cells.Item(48, 4).Value = "'0.003220"
$ws.Cells.Item(48, 5).Value = "'-9.43%"
$ws.Cells.Item(49, 4).Value = "'0.002280"
$ws.Cells.Item(49, 5).Value = "'-0.22%"
$ws.Cells.Item(50, 4).Value = "'0.00002102"
$ws.Cells.Item(50, 5).Value = "'-0.11%"
$ws.Cells.Item(51, 4).Value = "'0.0002002"
$ws.Cells.Item(51, 5).Value = "'-0.11%"
